$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "95.124.16"
$ws.Range("E2").Value = "  -1.69%  "
$ws.Range("D3").Value = "3.484.91"
$ws.Range("E3").Value = "  +4.36%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "646.18"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.46%  "
$ws.Range("E7").Value = "  +4.10%  "
$ws.Range("E8").Value = "  -3.84%  "
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.01"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.15%  "
$ws.Range("D11").Value = "3.480.62"
$ws.Range("E11").Value = "  +4.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.03"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.35%  "
$ws.Range("E13").Value = "  -3.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.20"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.60%  "
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "94.902.32"
$ws.Range("E15").Value = "  -1.50%  "
$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").Value = "4.144.74"
$ws.Range("E16").Value = "  +4.61%  "
$ws.Range("E17").Value = "  +2.14%  "
$ws.Range("E18").Value = "  -2.21%  "
$ws.Range("D19").Value = "3.479.65"
$ws.Range("E19").Value = "  +3.88%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.09"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.51"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +8.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.519"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -10.43%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "506.76"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.23%  "
$ws.Range("E24").Value = "  -4.45%  "
$ws.Range("E25").Value = "  -2.51%  "
$ws.Range("E26").Value = "  +1.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "95.47"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.15%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.20"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.73%  "
$ws.Range("D29").Value = "3.669.87"
$ws.Range("E29").Value = "  +4.41%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "11.92"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.01%  "
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.78"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +11.71%  "
$ws.Range("E33").Value = "  -4.86%  "
$ws.Range("E34").Value = "  -2.46%  "
$ws.Range("E35").Value = "  +9.72%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.573"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.04%  "
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "550.25"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.83%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.81"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.98%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.47"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.99%  "
$ws.Range("E41").Value = "  +13.37%  "
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("E43").Value = "  -0.30%  "
$ws.Range("E44").Value = "  -1.22%  "
$ws.Range("B45").Value = "ImmutableX"
$ws.Range("C45").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.72"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.65%  "
$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.69"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.51%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0418"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.37%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.54"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.79%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.18"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +9.97%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.24"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.73%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "53.42"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.05%  "
